$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add three new columns (precision / recall / F Measure) to both tables, plus
# a couple of corrected values in the second (averaged) table.
#
# The "NA" text must become shared-string index 15 (i.e. must be interned
# before "precision"/"recall"/"F Measure"), so populate the NA cells first.
# (PowerShell literal syntax here doesn't accept scientific notation, so all
# the small numbers below are written out in plain decimal form; they parse
# to the exact same IEEE-754 double as the scientific form.)
# ---------------------------------------------------------------------------

# --- NA placeholders (precision is undefined when there were no positive
#     predictions for that label) -------------------------------------------
$ws.Range("I3").Value = "NA"
$ws.Range("I6").Value = "NA"
$ws.Range("I7").Value = "NA"
$ws.Range("I8").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("I16").Value = "NA"
$ws.Range("I18").Value = "NA"

# --- Header row 1 (table #1) ------------------------------------------------
$ws.Range("I1").Value = "precision"
$ws.Range("J1").Value = "recall"
$ws.Range("K1").Value = "F Measure"
$ws.Range("I1:K1").Font.Bold = $true

# --- Header row 11 (table #2) -----------------------------------------------
$ws.Range("I11").Value = "precision"
$ws.Range("J11").Value = "recall"
$ws.Range("K11").Value = "F Measure"
$ws.Range("I11:K11").Font.Bold = $true

# --- Table #1 data rows (2-8): precision / recall / F Measure --------------
$ws.Range("I2").Value = 0.5
$ws.Range("J2").Value = 0.02272727
$ws.Range("K2").Value = 0.04347826

$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("I4").Value = 0.7675
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0

$ws.Range("I5").Value = 0.625
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0.86845827

$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.76923077

$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

# --- New row 9: a lone F Measure value of 0 ---------------------------------
$ws.Range("K9").Value = 0

# --- Table #2 corrected values (rows 12-18) ---------------------------------
$ws.Range("D12").Value = 0.600874844092425

$ws.Range("D13").Value = 0.603484231512754
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 1755
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = 624

$ws.Range("D14").Value = 0.559471917261171

$ws.Range("B15").Value = 0.640416666666667
$ws.Range("C15").Value = 0.359583333333333
$ws.Range("D15").Value = 0.522219057436343
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 863

$ws.Range("B16").Value = 0.757083333333333
$ws.Range("C16").Value = 0.242916666666667
$ws.Range("D16").Value = 0.504271172488533
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1817
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 583

$ws.Range("B17").Value = 0.649583333333333
$ws.Range("C17").Value = 0.350416666666667
$ws.Range("D17").Value = 0.545690247252748
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 1553
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 834

$ws.Range("D18").Value = 0.555486887418954

# --- Table #2 data rows (12-18): precision / recall / F Measure ------------
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0

$ws.Range("I13").Value = 0.619047619047619
$ws.Range("J13").Value = 0.0204081632653061
$ws.Range("K13").Value = 0.0395136778115502

$ws.Range("I14").Value = 0.711666666666667
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0.831548198636806

$ws.Range("I15").Value = 0.640266777824093
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 0.780686149936468

$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

$ws.Range("I17").Value = 0.461538461538462
$ws.Range("J17").Value = 0.00714285714285714
$ws.Range("K17").Value = 0.0140679953106682

$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0

# --- Cursor moves to K2 (matches the author's final selection) -------------
$ws.Range("K2").Select()
